$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add count column (C) values for the 3 rows that already have pin data
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 3

# Mark rows 13 and 14 as removed from the major cycle schedule ("XX")
$ws.Range("C13").Value = "XX"
$ws.Range("C14").Value = "XX"

# New rows 15 and 16 for newly added peripherals/comments
$ws.Range("C15").Value = "PWM0"
$ws.Range("A15").Value = "Pulse input"

$ws.Range("A16").Value = "Timing output"
$ws.Range("D16").Value = "PF3"
$ws.Range("E16").Value = "LED0"

$ws.Range("F35").Select()
